$d = $word.ActiveDocument

$d.Content.Find.Execute("Alma Nayeli Rodríguez Vázquez", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Clara Margarita Fernández Riveron", 2)
